# Circle Language Spec: Commands
# Merge "Command Object Referrers" and "Command Definition Referrers"
# into a single article.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1. Table "tblLook" cosmetic refresh (best effort; not independently
#    addressable through the exposed Table object model, so we just
#    re-assert the already-correct banding flags).
# ---------------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.ApplyStyleHeadingRows = $true
$tbl.ApplyStyleLastRow = $false
$tbl.ApplyStyleFirstColumn = $true
$tbl.ApplyStyleLastColumn = $false
$tbl.ApplyStyleRowBands = $true
$tbl.ApplyStyleColumnBands = $true

# ---------------------------------------------------------------------
# 2. Remove the stray "_GoBack" bookmark around the title run, by
#    replacing the title paragraph's XML without the bookmark tags.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleXml = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Circle Language Spec: Commands</w:t></w:r></w:p>'
$titlePara.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 3. Split the "Command Object Referrers" Heading2 paragraph into:
#      - Heading2 "Command Referrers"
#      - Heading3 "Command Object Referrers" (carrying the _GoBack
#        bookmark at the end of the paragraph)
# ---------------------------------------------------------------------
$headingPara = $d.Paragraphs(3)
$splitXml = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Command Referrers</w:t></w:r></w:p>' + `
            '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Command Object Referrers</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$headingPara.Range.InsertXML($splitXml)

# ---------------------------------------------------------------------
# 4. Append the new "Command Definition Referrers" article content at
#    the end of the body (before the final empty paragraph we also add).
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)

$newContent =
  '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Command Definition Referrers</w:t></w:r></w:p>' + `
  '<w:p ' + $wns + '><w:r><w:t>&lt; This topic will be further worked out in the future. &gt;</w:t></w:r></w:p>' + `
  '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Spacing"/></w:pPr></w:p>' + `
  '<w:p ' + $wns + '><w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Referrers </w:t></w:r><w:r><w:t xml:space="preserve">article explained how an object can be made aware of its referrers. A command is an object as well and the </w:t></w:r><w:r><w:rPr><w:rStyle w:val="CodeChar"/></w:rPr><w:t xml:space="preserve">Referrers </w:t></w:r><w:r><w:t xml:space="preserve">concept already provides a command with </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">command definition referrers </w:t></w:r><w:r><w:t xml:space="preserve">functionality. An inactive command is the only type of command symbol that can be referenced. Am inactive command symbol can implement the </w:t></w:r><w:r><w:rPr><w:rStyle w:val="CodeChar"/></w:rPr><w:t xml:space="preserve">Referrers </w:t></w:r><w:r><w:t>concept to register every call or reference to it.</w:t></w:r></w:p>' + `
  '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Spacing"/></w:pPr></w:p>' + `
  '<w:p ' + $wns + '><w:r><w:t xml:space="preserve">When a site hosts a command definition, that is widely used all over the world, you might not want the command definition to register its referrers, because it would be a very long list to maintain. You can turn off the </w:t></w:r><w:r><w:rPr><w:rStyle w:val="CodeChar"/></w:rPr><w:t xml:space="preserve">Referrers </w:t></w:r><w:r><w:t>concept for any command definition.</w:t></w:r></w:p>' + `
  '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Spacing"/></w:pPr></w:p>' + `
  '<w:p ' + $wns + '><w:r><w:t>If another site uses this widely used command definition, the using site could add a command reference to the command definition on the other site. A command reference has its own list of referrers. The using site could then redirect calls and references to its own command reference. Then the using site has a registration of anything on its site that uses the external command definition.</w:t></w:r></w:p>' + `
  '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>In a Diagram</w:t></w:r></w:p>' + `
  '<w:p ' + $wns + '><w:r><w:t xml:space="preserve">&lt; The expression of referrers in a diagram needs to be redone, because the referrers list refers to the parents of the references, which is not </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neccesarily</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the way to go. I' + [char]0x2019 + 'm not sure yet. I might want to register the related items and related lists items that are the references to the command definition, instead of registering their parents, and an ID, that the reference has inside the parent. &gt;</w:t></w:r></w:p>'

$newPara.Range.InsertXML($newContent)

# ---------------------------------------------------------------------
# 5. Trailing empty paragraph at the very end of the body.
# ---------------------------------------------------------------------
$tailPara = $d.Paragraphs($d.Paragraphs.Count)
$tailPara.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 6. Style metadata tweaks on the built-in "hidden" styles.
# ---------------------------------------------------------------------
$styleInfo = @(
    @{ Name = "Default Paragraph Font"; Priority = 1 },
    @{ Name = "Normal Table"; Priority = 99 },
    @{ Name = "No List"; Priority = 99 }
)
foreach ($info in $styleInfo) {
    $s = $d.Styles($info.Name)
    $s.Priority = $info.Priority
    $s.UnhideWhenUsed = $true
    try {
        $s.SemiHidden = $true
    } catch {
        # Not exposed by this COM host; best effort only.
    }
}
